$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the values in A2:A5 to the new set of IDs
$ws.Range("A2").Value = "A15011005700032C"
$ws.Range("A3").Value = "A15011005700025M"
$ws.Range("A4").Value = "A15011005800042C"
$ws.Range("A5").Value = "A15011005900012M"

# Remove the now-unused row 6 (previously A15011702300049)
$ws.Range("A6").Value = $null

# Update selection / dimension-related UI state to match new extent A1:A5
$ws.Range("A1:A5").Select()
